$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: group number
$ws.Range("D2").Value = 10

# "Eu" (me) - number and name
$ws.Range("D4").Value = 2182185
$ws.Range("G4").Value = "João Pedro Da Rocha Valverde Martins"

# 2º Elemento - number and name
$ws.Range("D6").Value = 2201793
$ws.Range("G6").Value = " Tiago José Figueira Pires Rodrigues dos Reis"

# 3º Elemento - number and name
$ws.Range("D7").Value = 2201790
$ws.Range("G7").Value = "Daniel Marques Gonçalves"

# Grupo de Funcionalidades do Projeto - implementation status
$ws.Range("I33").Value = "Parcial"
$ws.Range("I35").Value = "Parcial"
$ws.Range("I31").Value = "Completo"
$ws.Range("I32").Value = "Completo"
$ws.Range("I37").Value = "Completo"
